$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 112243462
$ws.Range("B2").Value = 78242
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"

$ws.Range("A3").Value = 112243461
$ws.Range("B3").Value = 77403
$ws.Range("E3").Value = 228912
$ws.Range("F3").Value = "Mörk kolflarnlav"
$ws.Range("G3").Value = "Carbonicola myrmecina"
$ws.Range("H3").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 410598
$ws.Range("R3").Value = 6710899

$ws.Range("A4").Value = 112243463
$ws.Range("B4").Value = 77403
$ws.Range("E4").Value = 228912
$ws.Range("F4").Value = "Mörk kolflarnlav"
$ws.Range("G4").Value = "Carbonicola myrmecina"
$ws.Range("H4").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q4").Value = 410608
$ws.Range("R4").Value = 6710914

$ws.Range("A5").Value = 112243468
$ws.Range("B5").Value = 77650
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 410566
$ws.Range("R5").Value = 6710872

$ws.Range("A6").Value = 112243469
$ws.Range("B6").Value = 77650
$ws.Range("Q6").Value = 410486
$ws.Range("R6").Value = 6710828

$ws.Range("A7").Value = 112243476
$ws.Range("B7").Value = 77403
$ws.Range("Q7").Value = 410524
$ws.Range("R7").Value = 6710795

$ws.Range("B8").Value = 78242
